$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 2501.889
$ws.Range("I55").Value = 359.4
$ws.Range("K55").Value = 359.4
$ws.Range("M55").Value = -145.4
$ws.Range("H68").Value = 111000
$ws.Range("J68").Value = 111000
$ws.Range("L68").Value = 111000
$ws.Range("N68").Value = -112498
$ws.Range("H71").Value = 111000
$ws.Range("J71").Value = 111000
$ws.Range("L71").Value = 333000
$ws.Range("N71").Value = -340488
$ws.Range("H80").Value = 1120.4166
$ws.Range("I80").Value = 1412.25
$ws.Range("J80").Value = 974.5
$ws.Range("K80").Value = 4236.75
$ws.Range("L80").Value = 2923.5
$ws.Range("M80").Value = -3238.75
$ws.Range("N80").Value = -4919.5
$ws.Range("H83").Value = 1120.4166
$ws.Range("I83").Value = 1412.25
$ws.Range("J83").Value = 974.5
$ws.Range("K83").Value = 12710.25
$ws.Range("L83").Value = 8770.5
$ws.Range("M83").Value = -7718.25
$ws.Range("N83").Value = -18754.5
$ws.Range("H116").Value = 5111
$ws.Range("I116").Value = 4953.7393
$ws.Range("J116").Value = 6316.6665
$ws.Range("K116").Value = 4953.7393
$ws.Range("L116").Value = 6316.6665
$ws.Range("M116").Value = -1511.7393
$ws.Range("N116").Value = -13200.6665
$ws.Range("H132").Value = 3302.5806
$ws.Range("I132").Value = 3195.9285
$ws.Range("K132").Value = 9587.7855
$ws.Range("M132").Value = -7057.7855
$ws.Range("H137").Value = 4938.0713
$ws.Range("I137").Value = 3085.8948
$ws.Range("K137").Value = 9257.6844
$ws.Range("M137").Value = -6707.6844
$ws.Range("H138").Value = 3229.9426
$ws.Range("J138").Value = 3348.2278
$ws.Range("L138").Value = 10044.6834
$ws.Range("N138").Value = -20324.6834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21795086
$ws.Range("I61").Value = 50007930
$ws.Range("K61").Value = 50007930
$ws.Range("M61").Value = -50007718
$ws.Range("H103").Value = 85000
$ws.Range("J103").Value = 85000
$ws.Range("L103").Value = 85000
$ws.Range("N103").Value = -87344
$ws.Range("H136").Value = 21795086
$ws.Range("I136").Value = 50007930
$ws.Range("K136").Value = 150023790
$ws.Range("M136").Value = -150021240

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 800.3
$ws.Range("I22").Value = 906.625
$ws.Range("K22").Value = 906.625
$ws.Range("M22").Value = -733.625
$ws.Range("H86").Value = 2267.25
$ws.Range("I86").Value = 1552.3077
$ws.Range("J86").Value = 3112.182
$ws.Range("K86").Value = 1552.3077
$ws.Range("L86").Value = 3112.182
$ws.Range("M86").Value = -429.3077000000001
$ws.Range("N86").Value = -5358.182
$ws.Range("H89").Value = 2267.25
$ws.Range("I89").Value = 1552.3077
$ws.Range("J89").Value = 3112.182
$ws.Range("K89").Value = 7761.538500000001
$ws.Range("L89").Value = 15560.91
$ws.Range("M89").Value = -2145.538500000001
$ws.Range("N89").Value = -26792.91
$ws.Range("H134").Value = 557569.9
$ws.Range("I134").Value = 1028.6
$ws.Range("K134").Value = 3085.8
$ws.Range("M134").Value = -550.7999999999997
$ws.Range("H135").Value = 60296.363
$ws.Range("J135").Value = 60296.363
$ws.Range("L135").Value = 60296.363
$ws.Range("N135").Value = -70436.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8486.25
$ws.Range("I99").Value = 7975
$ws.Range("K99").Value = 7975
$ws.Range("M99").Value = -6477
$ws.Range("H112").Value = 87516.164
$ws.Range("J112").Value = 87516.164
$ws.Range("L112").Value = 87516.164
$ws.Range("N112").Value = -90470.164
$ws.Range("H126").Value = 8486.25
$ws.Range("I126").Value = 7975
$ws.Range("K126").Value = 23925
$ws.Range("M126").Value = -21455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 112.63415
$ws.Range("I2").Value = 83.47369
$ws.Range("J2").Value = 137.81818
$ws.Range("K2").Value = 500.84214
$ws.Range("L2").Value = 826.9090800000001
$ws.Range("M2").Value = -387.84214
$ws.Range("N2").Value = -1052.90908
$ws.Range("H10").Value = 419.13635
$ws.Range("I10").Value = 40
$ws.Range("J10").Value = 479
$ws.Range("K10").Value = 120
$ws.Range("L10").Value = 1437
$ws.Range("M10").Value = 19
$ws.Range("N10").Value = -1715
$ws.Range("H88").Value = 3808.3333
$ws.Range("J88").Value = 3990
$ws.Range("L88").Value = 11970
$ws.Range("N88").Value = -12826
$ws.Range("H91").Value = 3808.3333
$ws.Range("J91").Value = 3990
$ws.Range("L91").Value = 11970
$ws.Range("N91").Value = -14934
$ws.Range("H131").Value = 5812.488
$ws.Range("I131").Value = 6812.25
$ws.Range("J131").Value = 5570.121
$ws.Range("K131").Value = 20436.75
$ws.Range("L131").Value = 16710.363
$ws.Range("M131").Value = -15396.75
$ws.Range("N131").Value = -26790.363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 7687.5
$ws.Range("J28").Value = 8583.333000000001
$ws.Range("L28").Value = 8583.333000000001
$ws.Range("N28").Value = -8967.333000000001
$ws.Range("H112").Value = 98500
$ws.Range("J112").Value = 98500
$ws.Range("L112").Value = 98500
$ws.Range("N112").Value = -100716
$ws.Range("H122").Value = 1214.2609
$ws.Range("I122").Value = 978.8333
$ws.Range("J122").Value = 2061.8
$ws.Range("K122").Value = 2936.4999
$ws.Range("L122").Value = 6185.400000000001
$ws.Range("M122").Value = -486.4998999999998
$ws.Range("N122").Value = -11085.4
$ws.Range("H132").Value = 100016120
$ws.Range("I132").Value = 250003100
$ws.Range("J132").Value = 24801.834
$ws.Range("K132").Value = 750009300
$ws.Range("L132").Value = 74405.50199999999
$ws.Range("M132").Value = -750006770
$ws.Range("N132").Value = -79465.50199999999
$ws.Range("H140").Value = 82779.39999999999
$ws.Range("J140").Value = 82779.39999999999
$ws.Range("L140").Value = 82779.39999999999
$ws.Range("N140").Value = -93139.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1754.1538
$ws.Range("I61").Value = 1484.6666
$ws.Range("K61").Value = 1484.6666
$ws.Range("M61").Value = -1282.6666
$ws.Range("H113").Value = 1754.1538
$ws.Range("I113").Value = 1484.6666
$ws.Range("K113").Value = 1484.6666
$ws.Range("M113").Value = 685.3334
$ws.Range("H122").Value = 5863.25
$ws.Range("I122").Value = 5642.1113
$ws.Range("K122").Value = 16926.3339
$ws.Range("M122").Value = -14476.3339
$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 3416657.8
$ws.Range("J75").Value = 99989.39999999999
$ws.Range("L75").Value = 99989.39999999999
$ws.Range("N75").Value = -101861.4
$ws.Range("H78").Value = 3416657.8
$ws.Range("J78").Value = 99989.39999999999
$ws.Range("L78").Value = 299968.2
$ws.Range("N78").Value = -309328.2
$ws.Range("H114").Value = 90988.336
$ws.Range("J114").Value = 90988.336
$ws.Range("L114").Value = 90988.336
$ws.Range("N114").Value = -99666.336
$ws.Range("H122").Value = 3612.8
$ws.Range("I122").Value = 3442.2856
$ws.Range("K122").Value = 10326.8568
$ws.Range("M122").Value = -7876.856800000001
